# "Generate Report for Archive"
#
# 1) The handoff status text changed from "Ready for handoff" to
#    "In Translation" everywhere it is used: the Overview sheet's
#    per-language status columns (E2 = zh-cn, F2 = de-de) and each
#    language sheet's own Status column (C2).
#
# 2) Because the new status text is shorter than the old one, the
#    Status column(s) were narrowed accordingly:
#      - Overview!E:E and Overview!F:F
#      - zh-cn!C:C
#      - de-de!C:C
#    all shrank from ~17.22 characters wide to ~13.41 characters wide.

$wb = $excel.ActiveWorkbook

$newStatus = "In Translation"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Update the status text ------------------------------------------------
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsZhCn.Range("C2").Value = $newStatus
$wsDeDe.Range("C2").Value = $newStatus

# --- Narrow the Status columns to match the shorter text --------------------
# Target stored width is ~13.41 characters; the COM ColumnWidth property
# here only resolves to the nearest 1/6 of a character, so 12.5 (-> ~13.33
# stored) is the closest attainable value to that target.
$newColumnWidth = 12.5

$wsOverview.Columns.Item(5).ColumnWidth = $newColumnWidth  # column E
$wsOverview.Columns.Item(6).ColumnWidth = $newColumnWidth  # column F
$wsZhCn.Columns.Item(3).ColumnWidth = $newColumnWidth      # column C
$wsDeDe.Columns.Item(3).ColumnWidth = $newColumnWidth      # column C
